{"js": "// Update the p-value tables in \"Fig 8\" (Cod = table 0, Hake = table 1).\n// Each entry is [tableIndex, rowIndex, columnIndex, expectedOldText, newText].\n// Row/column indices are 0-based and include the header row/column, matching\n// Word.Table.getCell()'s own (row, column) addressing.\nconst edits = [\n  [0, 1, 1, \"0.06\", \"0.10\"],\n  [0, 1, 3, \"0.03\", \"0.07\"],\n  [0, 2, 1, \"0.71\", \"0.75\"],\n  [0, 2, 3, \"0.85\", \"0.97\"],\n  [0, 3, 1, \"0.97\", \"0.51\"],\n  [0, 3, 3, \"0.47\", \"0.34\"],\n  [0, 4, 1, \"0.81\", \"0.77\"],\n  [0, 4, 3, \"0.03\", \"0.06\"],\n  [0, 5, 1, \"0.77\", \"0.69\"],\n  [0, 5, 3, \"0.02\", \"0.03\"],\n  [1, 1, 3, \"0.84\", \"0.94\"],\n  [1, 2, 3, \"0.62\", \"0.52\"],\n  [1, 3, 3, \"0.37\", \"0.25\"],\n  [1, 4, 3, \"0.75\", \"0.64\"],\n  [1, 5, 3, \"0.12\", \"0.10\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Grab every target cell and confirm its current text matches what the diff\n// expects to find there (guards against drift in table layout breaking the\n// positional mapping), then replace just the text of that run via a\n// cell-scoped search-and-replace so the existing run formatting (font,\n// size, color, xml:space) is preserved instead of being rebuilt from\n// scratch.\nfor (const [tableIndex, rowIndex, colIndex, oldText, newText] of edits) {\n  const cell = tables.items[tableIndex].getCell(rowIndex, colIndex);\n  cell.load(\"value\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n\n  if (cell.value !== oldText) {\n    throw new Error(\n      `Cell mismatch at table ${tableIndex} row ${rowIndex} col ${colIndex}: ` +\n        `expected \"${oldText}\" but found \"${cell.value}\"`\n    );\n  }\n\n  const hits = cell.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  hits.load(\"items\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n\n  if (hits.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${oldText}\" in table ${tableIndex} row ${rowIndex} col ${colIndex}, found ${hits.items.length}`\n    );\n  }\n\n  hits.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the p-value tables in \"Fig 8\" (Table 1 = Cod, Table 2 = Hake).\n# Each entry is (tableIndex, rowIndex, columnIndex, expectedOldText, newText)\n# using Word's 1-based Tables/Cell indexing (row 1 = header row).\n$d = $word.ActiveDocument\n\n$edits = @(\n  @(1, 2, 2, \"0.06\", \"0.10\"),\n  @(1, 2, 4, \"0.03\", \"0.07\"),\n  @(1, 3, 2, \"0.71\", \"0.75\"),\n  @(1, 3, 4, \"0.85\", \"0.97\"),\n  @(1, 4, 2, \"0.97\", \"0.51\"),\n  @(1, 4, 4, \"0.47\", \"0.34\"),\n  @(1, 5, 2, \"0.81\", \"0.77\"),\n  @(1, 5, 4, \"0.03\", \"0.06\"),\n  @(1, 6, 2, \"0.77\", \"0.69\"),\n  @(1, 6, 4, \"0.02\", \"0.03\"),\n  @(2, 2, 4, \"0.84\", \"0.94\"),\n  @(2, 3, 4, \"0.62\", \"0.52\"),\n  @(2, 4, 4, \"0.37\", \"0.25\"),\n  @(2, 5, 4, \"0.75\", \"0.64\"),\n  @(2, 6, 4, \"0.12\", \"0.10\")\n)\n\nforeach ($edit in $edits) {\n  $tableIndex = $edit[0]\n  $rowIndex = $edit[1]\n  $colIndex = $edit[2]\n  $oldText = $edit[3]\n  $newText = $edit[4]\n\n  $cell = $d.Tables.Item($tableIndex).Cell($rowIndex, $colIndex)\n  $range = $cell.Range\n  # Cell.Range.Text includes the trailing cell-mark (CR + BEL); strip it\n  # before comparing so we can confirm we are about to overwrite the value\n  # the diff expects, guarding against drift in the table layout.\n  $currentText = $range.Text.TrimEnd([char]13, [char]7)\n\n  if ($currentText -ne $oldText) {\n    throw \"Cell mismatch at table $tableIndex row $rowIndex col $colIndex`: expected '$oldText' but found '$currentText'\"\n  }\n\n  $range.Text = $newText\n}\n"}
